$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Structural edits: insert a new column before C (for "scissor"),
#    and insert two new rows before row 18 (new "all" + scissor-only
#    sample rows). Excel's COM Insert() shifts existing cell content,
#    formulas and column widths automatically.
# -----------------------------------------------------------------
$ws.Columns("C:C").Insert()
$ws.Rows("18:19").Insert()

# -----------------------------------------------------------------
# 2. New header for the inserted "scissor" column.
# -----------------------------------------------------------------
$ws.Range("C1").Value = "scissor"

# -----------------------------------------------------------------
# 3. New row 18 - "all" columns checked (maj7, scissor, leveller,
#    maj7width, echo, cathedral, smasher all = x), plus sizes.
# -----------------------------------------------------------------
$ws.Range("A18").Value = "all"
$ws.Range("B18").Value = "x"
$ws.Range("C18").Value = "x"
$ws.Range("D18").Value = "x"
$ws.Range("E18").Value = "x"
$ws.Range("F18").Value = "x"
$ws.Range("G18").Value = "x"
$ws.Range("H18").Value = "x"
$ws.Range("I18").Value = 48176
$ws.Range("J18").Value = 20036
$ws.Range("J18").Interior.Color = 65535

# -----------------------------------------------------------------
# 4. New row 19 - scissor-only sample row.
# -----------------------------------------------------------------
$ws.Range("C19").Value = "x"
$ws.Range("I19").Value = 23192
$ws.Range("J19").Value = 8564

# -----------------------------------------------------------------
# 5. New rows 23 and 24 (filling in previously blank rows between the
#    "filter static letters" row, now 22, and the summary table).
# -----------------------------------------------------------------
$ws.Range("A23").Value = "static init of maj7width"
$ws.Range("J23").Value = 19680
$ws.Range("K23").Formula = '=$J$12-J23'

$ws.Range("A24").Value = "adding features to maj7width-"
$ws.Range("J24").Value = 19728
$ws.Range("K24").Formula = '=$J$12-J24'

# -----------------------------------------------------------------
# 6. New row 38 in the summary table, for "scissor".
# -----------------------------------------------------------------
$ws.Range("F38").Value = "scissor"
$ws.Range("G38").Formula = '=J19-J2'
$ws.Range("H38").Formula = '=J18-J12'

# -----------------------------------------------------------------
# 7. Fix up conditional-formatting ranges that don't auto-shift with
#    the column insert (engine leaves the legacy AppliesTo range
#    pointing at the pre-insert columns).
# -----------------------------------------------------------------
$cf1 = $ws.Range("H1:I1048576").FormatConditions.Item(1)
$cf1.ModifyAppliesToRange($ws.Range("I1:J1048576"))

$cf2 = $ws.Range("J1:J1048576").FormatConditions.Item(1)
$cf2.ModifyAppliesToRange($ws.Range("K1:K1048576"))

# -----------------------------------------------------------------
# 8. Restore the view's active selection to match the author's.
# -----------------------------------------------------------------
$ws.Range("J16").Select()

Write-Output "edit complete"
